$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ORLY")

# Row 4 - Inventory
$ws.Range("B4").Value = 3622000000.0
$ws.Range("C4").Value = 3653000000.0
$ws.Range("D4").Value = 3527000000.0
$ws.Range("E4").Value = 3529000000.0
$ws.Range("F4").Value = 3557000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 4318000000.0
$ws.Range("C14").Value = 4185000000.0
$ws.Range("D14").Value = 4084000000.0
$ws.Range("E14").Value = 3936000000.0
$ws.Range("F14").Value = 3758000000.0

# Row 23 - Long Term Tax Liability (Deferred)
$ws.Range("B23").Value = 165000000.0
$ws.Range("C23").Value = 156000000.0
$ws.Range("D23").Value = 174000000.0
$ws.Range("E23").Value = 156000000.0
$ws.Range("F23").Value = 73000000.0

# Row 35 - Net Debt
$ws.Range("G35").Value = 4166121000.0

# Row 36 - Total Debt
$ws.Range("G36").Value = 4206527000.0
